# update 22 agustus 2023
# Adds new "Kas HIMA" ledger rows (28 Mei 2023 .. 18 Agustus 2023) into Sheet1,
# pushing the previous last data row (6 Juni 2023) down by one row and filling
# in the blank placeholder rows that already existed below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 26 first gets what used to live in row 25 (6 Jun 2023 / 39000 masuk),
# then row 25 becomes the newly-recorded 28 Mei 2023 entry (7500 keluar -
# biaya admin rekening). Writing row 26 before row 25 means we never lose the
# original values while we "shift" them down.
# ---------------------------------------------------------------------------

# Row 26 (was row 25): 6 Juni 2023, masuk 39000, uang persembahan - reguler
$ws.Range("A26").Value2 = 45083
$ws.Range("B26").Value2 = 39000
$ws.Range("B26").NumberFormat = $ws.Range("B25").NumberFormat
$ws.Range("C26").Value2 = 0
$ws.Range("D26").Formula = "=D25+B26-C26"
$ws.Range("E26").Value2 = "uang persembahan - reguler"
$ws.Range("F26").Value2 = "yofandi"

# Row 25 (new): 28 Mei 2023, keluar 7500, biaya admin rekening
$ws.Range("A25").Value2 = 45074
$ws.Range("B25").Value2 = 0
$ws.Range("C25").Value2 = 7500
$ws.Range("D25").Formula = "=D24+B25-C25"
$ws.Range("E25").Value2 = "biaya admin rekening"
$ws.Range("F25").Value2 = "yofandi"

# Row 27 (new, fills existing blank placeholder row): 16 Juni 2023, keluar 7500
$ws.Range("A27").Value2 = 45093
$ws.Range("C27").Value2 = 7500
$ws.Range("D27").Formula = "=D26+B27-C27"
$ws.Range("E27").Value2 = "biaya admin rekening"
$ws.Range("F27").Value2 = "yofandi"

# Row 28 (new): 19 Juni 2023, masuk 73000, uang persembahan - reguler
$ws.Range("A28").Value2 = 45096
$ws.Range("B28").Value2 = 73000
$ws.Range("B28").NumberFormat = $ws.Range("B25").NumberFormat
$ws.Range("C28").Value2 = 0
$ws.Range("D28").Formula = "=D27+B28-C28"
$ws.Range("E28").Value2 = "uang persembahan - reguler"
$ws.Range("F28").Value2 = "yofandi"

# Row 29 (new): 3 Juli 2023, masuk 30000, uang persembahan - reguler
$ws.Range("A29").Value2 = 45110
$ws.Range("B29").Value2 = 30000
$ws.Range("B29").NumberFormat = $ws.Range("B25").NumberFormat
$ws.Range("C29").Value2 = 0
$ws.Range("D29").Formula = "=D28+B29-C29"
$ws.Range("E29").Value2 = "uang persembahan - reguler"
$ws.Range("F29").Value2 = "yofandi"

# Row 30 (new): 17 Juli 2023, masuk 30000, uang persembahan - reguler
$ws.Range("A30").Value2 = 45124
$ws.Range("B30").Value2 = 30000
$ws.Range("B30").NumberFormat = $ws.Range("B25").NumberFormat
$ws.Range("C30").Value2 = 0
$ws.Range("D30").Formula = "=D29+B30-C30"
$ws.Range("E30").Value2 = "uang persembahan - reguler"
$ws.Range("F30").Value2 = "yofandi"

# Row 31 (new): 22 Juli 2023, keluar 7500, biaya admin rekening
$ws.Range("A31").Value2 = 45129
$ws.Range("B31").Value2 = 0
$ws.Range("B31").NumberFormat = $ws.Range("B25").NumberFormat
$ws.Range("C31").Value2 = 7500
$ws.Range("D31").Formula = "=D30+B31-C31"
$ws.Range("E31").Value2 = "biaya admin rekening"
$ws.Range("F31").Value2 = "yofandi"

# Row 32 (new): 18 Agustus 2023, keluar 7500, biaya admin rekening
$ws.Range("A32").Value2 = 45156
$ws.Range("B32").Value2 = 0
$ws.Range("B32").NumberFormat = $ws.Range("B25").NumberFormat
$ws.Range("C32").Value2 = 7500
$ws.Range("D32").Formula = "=D31+B32-C32"
$ws.Range("E32").Value2 = "biaya admin rekening"
$ws.Range("F32").Value2 = "yofandi"

# Keep one more trailing blank styled row at the bottom of the sheet (row 34
# -> 35), matching the rest of the blank "C" column placeholders.
$ws.Rows(34).Copy()
$ws.Rows(34).Insert(-4121) # xlShiftDown

# Move the view roughly to where the edits were made.
$ws.Range("H32").Select()
